$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# Edit the employee record so its CompanyName matches the existing
# "ActivityCompany" employee record instead of "StandardTestCompany".
$ws.Range("A2").Value = "ActivityCompany"
$ws.Range("A3").Value = "ActivityCompany"

# Match the author's final cursor position/selection on the sheet.
$ws.Activate()
[void]$ws.Range("B8").Select()
